$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# --- Row 21: update the remarks cell (cursor-control stick FPC spec) ---
$ws.Range("G21").Value = "カーソル操作用 / FPC 4ピン 1.0mmピッチ / VCC=3.3V直結OK / ピンアサイン: GND-X-VCC-Y / 別途購入（AliExpress or 中古）"

# --- Row 22: update the FPC connector part (6P/0.5mm -> 4P/1.0mm) ---
$ws.Range("C22").Value = "Molex 5034800440"

# D22 is the purely-numeric-looking part number "5034800440". A plain
# .Value assignment would be auto-coerced by Excel into a numeric cell,
# but the source part number must remain text (matching the original
# cell's text type). Write it as a text-producing formula first, then
# flatten the formula down to a static value via copy/paste-special so
# the stored cell is plain text, not a live formula.
$ws.Range("D22").Formula = '="5034800440"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("E22").Value = "C3170007"
$ws.Range("G22").Value = "4P 1.0mmピッチ ZIF ヒンジ式 両面接触 / 3DSスライドパッド接続用 / JLCPCB在庫928個"
